$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$newRow = 18

# Column A ("2025-01-10") looks like a date and column D ("01") looks like a
# number, so Excel would otherwise auto-convert them. Force them to be
# entered as plain text, then reset the cell style back to Normal so no
# extra number-format style sticks to the cell (matching the plain,
# unstyled text cells used elsewhere in this sheet).
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025-01-10"
$ws.Range("A$newRow").Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "22:30:36"
$ws.Cells.Item($newRow, 3).Value = "Friday"

$ws.Range("D$newRow").NumberFormat = "@"
$ws.Cells.Item($newRow, 4).Value = "01"
$ws.Range("D$newRow").Style = "Normal"

# Columns E-T are numeric values
$ws.Cells.Item($newRow, 5).Value = 127139
$ws.Cells.Item($newRow, 6).Value = 143706
$ws.Cells.Item($newRow, 7).Value = 169671
$ws.Cells.Item($newRow, 8).Value = 159722
$ws.Cells.Item($newRow, 9).Value = -1
$ws.Cells.Item($newRow, 10).Value = 142939
$ws.Cells.Item($newRow, 11).Value = -1
$ws.Cells.Item($newRow, 12).Value = -1
$ws.Cells.Item($newRow, 13).Value = 192926
$ws.Cells.Item($newRow, 14).Value = 115393
$ws.Cells.Item($newRow, 15).Value = 45859
$ws.Cells.Item($newRow, 16).Value = 28491
$ws.Cells.Item($newRow, 17).Value = 65238
$ws.Cells.Item($newRow, 18).Value = -1
$ws.Cells.Item($newRow, 19).Value = 49164
$ws.Cells.Item($newRow, 20).Value = -1
